$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The header in A1 was "URN:Tel" and should become "phone"
$ws.Range("A1").Value = "phone"
